$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the Accuracy/Loss value columns (D:E, rows 2-11) two columns to the
# right (F:G), leaving the header row (row 1) and the bottom style-only
# row (row 13) untouched.
$src = $ws.Range("D2:E11")
$dst = $ws.Range("F2")
$src.Cut($dst) | Out-Null

# Clear whatever leftover formatting/content Cut left behind in D2:E11.
$ws.Range("D2:E11").Clear() | Out-Null

# Match the author's final selection.
$ws.Range("D2").Select() | Out-Null
